$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 36) -----------------------------------------
# Row 33 (A33:I33) already carries the exact same cell-style combination that
# the new row needs (date/border/fill pattern used by the most recent
# entries), so copy its formatting down into row 36 before filling in values.
$ws.Range("A33:I33").Copy()
$ws.Range("A36:I36").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(36, 1).Value = 45046
$ws.Cells.Item(36, 2).Value = "BUAM"
$ws.Cells.Item(36, 3).Value = "N/A"
$ws.Cells.Item(36, 4).Value = "Sainte-Ours"
$ws.Cells.Item(36, 5).Value = "Montérégie"
$ws.Cells.Item(36, 6).Value = "A"
$ws.Cells.Item(36, 7).Value = "NA"
$ws.Cells.Item(36, 8).Value = "Amplexus et ponte"
$ws.Cells.Item(36, 9).Value = "Alexandre Gariépy/Marc DuBois (Facebook, Amphibiens et Reptiles du Québec)"

# --- Column width tweaks (G widened, I widened) -----------------------------
$ws.Columns.Item(7).ColumnWidth = 14.42578125
$ws.Columns.Item(9).ColumnWidth = 80.5703125

# --- Selection moves to H40 (as left in the source file) -------------------
[void]$ws.Range("H40").Select()
